# Adds KNN model save results (KNN - Red / KNN - White) to the "Model Review"
# sheet, mirroring the pattern already used for the other model rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model Review")

# ---- Row 6: KNN - Red ----------------------------------------------------
$ws.Range("C6").Value = "Standard"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.74375000000000002
$ws.Range("F6").Value = 0.74
$ws.Range("H6").Value = "Top Features: Alcohol, sulphates. Free sulfur dioxide, chlorides, volatile acidity, total sulfure dioxide`nBest Params: Best leaf_size: 1 Best p: 2 Best n_neighbors: 1"
$ws.Range("J6").Value = 42

# ---- Row 7: KNN - White ----------------------------------------------------
$ws.Range("C7").Value = "Standard"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.79110000000000003
$ws.Range("F7").Value = 0.79
$ws.Range("H7").Value = "Top Features: All with Chlorides having the least weight`nBest Params: n_neighbors = 1, leaf_size = 1, p = 1"
$ws.Range("J7").Value = 42

# The Notes column (H) wraps text; with real Excel this auto-sizes the row
# height.  Set the resulting heights explicitly so the rows render the same
# way as the other multi-line note rows (8 and 9) on this sheet.
$ws.Rows.Item(6).RowHeight = 86.4
$ws.Rows.Item(7).RowHeight = 72

# Restore the cursor/selection left behind after the edit.
[void]$ws.Range("D8").Select()
